# Applies numeric "want-to-go" count bumps (col F) and several status/
# listing updates to the 杭州-漫展信息 workbook, per the upstream data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12776
$ws.Range("F3").Value = 7140
$ws.Range("F6").Value = 448
$ws.Range("C7").Value = "杭州·少女番only3.0（取消）"
$ws.Range("G7").Value = "不可售"
$ws.Range("F12").Value = 352
$ws.Range("F13").Value = 1008
$ws.Range("F14").Value = 3
$ws.Range("F16").Value = 1017
$ws.Range("F18").Value = 242
$ws.Range("F19").Value = 368
$ws.Range("F21").Value = 276
$ws.Range("F22").Value = 308
$ws.Range("F24").Value = 151
$ws.Range("F25").Value = 369
$ws.Range("F26").Value = 5230
$ws.Range("F28").Value = 1424
$ws.Range("F29").Value = 307
$ws.Range("F30").Value = 1339
$ws.Range("F31").Value = 61
$ws.Range("F32").Value = 34
$ws.Range("F33").Value = 1356
$ws.Range("F35").Value = 4
$ws.Range("F38").Value = 3729
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 4
$ws.Range("F19").Value = 19
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9275
$ws.Range("F3").Value = 559
$ws.Range("F4").Value = 2004
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9275
$ws.Range("F3").Value = 559
$ws.Range("F4").Value = 2004
$ws.Range("F5").Value = 12776
$ws.Range("F6").Value = 7140
$ws.Range("F9").Value = 448
$ws.Range("C10").Value = "浙江·蔚蓝档案ONLY02-夏末狂欢！"
$ws.Range("D10").Value = "康候圣街99号 顺丰创新中心"
$ws.Range("F10").Value = 1000
$ws.Range("G10").Value = "已售罄"
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=86594"
$ws.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202405/TVvJFURG1716799911888.jpeg"
$ws.Range("F12").Value = 352
$ws.Range("F13").Value = 1008
$ws.Range("F14").Value = 3
$ws.Range("F16").Value = 1017
$ws.Range("F18").Value = 242
$ws.Range("F19").Value = 368
$ws.Range("F21").Value = 276
$ws.Range("F22").Value = 308
$ws.Range("F26").Value = 4
$ws.Range("F27").Value = 151
$ws.Range("F28").Value = 369
$ws.Range("F29").Value = 5230
$ws.Range("F31").Value = 1424
$ws.Range("F34").Value = 307
$ws.Range("F36").Value = 1339
$ws.Range("F37").Value = 61
$ws.Range("F38").Value = 1356
$ws.Range("F47").Value = 3729
$ws.Range("F48").Value = 19
